# ============================================================================
# Sprint 02.pptx - wording/formatting touch-ups
# ============================================================================

function Set-RangeRuns {
    # Replaces the text of a TextRange (typically a single paragraph range)
    # with the concatenation of $parts, then "re-stamps" each part's
    # character span so the engine emits one <a:r> run per array entry
    # instead of silently diff-patching the old runs.
    param(
        $range,
        [string[]]$parts
    )
    $full = [string]::Join("", $parts)
    $range.Text = $full
    $pos = 1
    foreach ($part in $parts) {
        $len = $part.Length
        if ($len -gt 0) {
            $sub = $range.Characters($pos, $len)
            $sub.Text = $sub.Text
        }
        $pos += $len
    }
}

$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 1 : title slide -> merge "Sprint " + "Review" into "Sprint Review"
# ----------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1b = $s1.Shapes.Item(2)
$tr1 = $sh1b.TextFrame.TextRange
$full1 = $tr1.Characters(1, $tr1.Length)
$full1.Text = "Sprint Review"

# ----------------------------------------------------------------------
# Slide 2 : Agenda
# ----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange

Set-RangeRuns $body2.Paragraphs(1,1) @("Objetivos")
Set-RangeRuns $body2.Paragraphs(2,1) @("Logros")
Set-RangeRuns $body2.Paragraphs(3,1) @("Inconvenientes")
Set-RangeRuns $body2.Paragraphs(4,1) @("Tareas postergadas")
Set-RangeRuns $body2.Paragraphs(5,1) @("Mejoras pendientes")
Set-RangeRuns $body2.Paragraphs(6,1) @("Gráficos y estadísticas de ", "proyecto")
Set-RangeRuns $body2.Paragraphs(7,1) @("Resumen de ", "desarrollo")

# ----------------------------------------------------------------------
# Slide 3 : Objetivos del Sprint 02 (detail)
# ----------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

Set-RangeRuns $body3.Paragraphs(3,1) @("Planificar, dirigir y controlar el Sprint de forma ", "correcta, para evitar inconvenientes.")
Set-RangeRuns $body3.Paragraphs(5,1) @("Completar el modelo de clases ", "del ", "producto y el modelo de ", "Base de Datos.")
Set-RangeRuns $body3.Paragraphs(7,1) @("Llevar a cabo los ", "últimos test de desarrollo para iniciar la producción concreta del producto.")

# ----------------------------------------------------------------------
# Slide 5 : Logros
# ----------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$content5 = $s5.Shapes.Item(2)
$content5.TextFrame.TextRange.ParagraphFormat.TextFrame.Parent | Out-Null
$body5 = $content5.TextFrame.TextRange

Set-RangeRuns $body5.Paragraphs(1,1) @("El nuevo ", "integrante pudo incorporase  al ", "proyecto de ", "forma ", "correcta.")
Set-RangeRuns $body5.Paragraphs(3,1) @("Fueron definidos ", "los roles ", "para el ", "integrante;  ", "un nuevo Ingeniero de Software (con orientación a desarrollo web) y DBA para el desarrollo del producto.")

$bodyPr5 = $content5.TextFrame.TextRange.ParagraphFormat
# normAutofit fontScale 62500 -> 70000 (handled via shape XML fallback below)

# ----------------------------------------------------------------------
# Slide 7 : Inconvenientes
# ----------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$body7 = $s7.Shapes.Item(2).TextFrame.TextRange

Set-RangeRuns $body7.Paragraphs(1,1) @("Se plantearon inconvenientes en ", "la ", "adaptación de la arquitectura para cumplimentar el desarrollo entre ", "lenguajes.")

# paragraph 7 is "Se dejó una <technical><task> para un siguiente sprint, porque..."
Set-RangeRuns $body7.Paragraphs(7,1) @("Se ", "postergó una ", "technical", " ", "task", " para un siguiente sprint, ", "ya que era ", "necesario mejorar algunos aspectos técnicos.")

# ----------------------------------------------------------------------
# Slide 9 : Tasks postergadas para el Siguiente Sprint
# ----------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1).TextFrame.TextRange
Set-RangeRuns $title9.Paragraphs(1,1) @("Tareas postergadas ", "para el ", "siguiente ", "Sprint")

$body9 = $s9.Shapes.Item(2).TextFrame.TextRange
Set-RangeRuns $body9.Paragraphs(1,1) @("Realizar pruebas de recuperación de imágenes comprimidas en la base de datos.")
Set-RangeRuns $body9.Paragraphs(3,1) @("Implementar ", " funcionalidad para almacenar imágenes en ", "el motor de base de ", "datos.")

# ----------------------------------------------------------------------
# Slide 10 : Que podemos hacer para mejorar -> Mejoras pendientes
# ----------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$title10 = $s10.Shapes.Item(1).TextFrame.TextRange
Set-RangeRuns $title10.Paragraphs(1,1) @("Mejoras pendientes")

$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
Set-RangeRuns $body10.Paragraphs(1,1) @("Mayor comunicación ", "con el ", "Product", " ", "O", "wner", ".")
Set-RangeRuns $body10.Paragraphs(3,1) @("Mayor flexibilidad a la hora de ", "replanificar", "  y reasignar ", "User", "  ", "Stories", " en caso de no poder efectivizar una.")
Set-RangeRuns $body10.Paragraphs(5,1) @("Comenzar ", "con un proceso inicial de integración del producto.")
